# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns.
# Values that would otherwise be auto-parsed by Excel as a plain number
# (single-decimal-point numerics) are written with a leading apostrophe
# so they are kept as text, matching the original inline-string data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "65.061.75";  E = "  -0.22%  " }
    @{ Row = 3;  D = "3.572.25";   E = "  +2.79%  " }
    @{ Row = 4;  D = $null;        E = "  -0.08%  " }
    @{ Row = 5;  D = "600.98";     E = "  +2.16%  " }
    @{ Row = 6;  D = "135.87";     E = "  -1.16%  " }
    @{ Row = 7;  D = "3.570.42";   E = "  +2.81%  " }
    @{ Row = 8;  D = $null;        E = "  -0.02%  " }
    @{ Row = 9;  D = $null;        E = "  +1.06%  " }
    @{ Row = 10; D = $null;        E = "  +0.81%  " }
    @{ Row = 11; D = "6.93";       E = "  -3.62%  " }
    @{ Row = 12; D = $null;        E = "  +1.60%  " }
    @{ Row = 13; D = "4.176.50";   E = "  +2.70%  " }
    @{ Row = 14; D = $null;        E = "  +0.42%  " }
    @{ Row = 15; D = "3.578.50";   E = "  +4.32%  " }
    @{ Row = 16; D = "27.16";      E = "  +2.10%  " }
    @{ Row = 17; D = $null;        E = "  +0.56%  " }
    @{ Row = 18; D = "65.135.05";  E = "  +0.07%  " }
    @{ Row = 19; D = "10.11";      E = "  +3.58%  " }
    @{ Row = 20; D = $null;        E = "  +4.03%  " }
    @{ Row = 21; D = $null;        E = "  +1.64%  " }
    @{ Row = 22; D = "388.86";     E = "  -0.23%  " }
    @{ Row = 23; D = $null;        E = "  +4.45%  " }
    @{ Row = 24; D = "3.717.44";   E = "  +2.75%  " }
    @{ Row = 25; D = "74.22";      E = "  +2.31%  " }
    @{ Row = 26; D = $null;        E = "  +0.04%  " }
    @{ Row = 27; D = "0.0000117";  E = "  +6.33%  " }
    @{ Row = 28; D = "7.75";       E = "  +4.89%  " }
    @{ Row = 29; D = "1.00";       E = "  +0.29%  " }
    @{ Row = 30; D = $null;        E = "  +3.23%  " }
    @{ Row = 31; D = "8.43";       E = "  +2.64%  " }
    @{ Row = 32; D = $null;        E = "  +22.05%  " }
    @{ Row = 33; D = "3.577.88";   E = "  +2.49%  " }
    @{ Row = 34; D = "24.03";      E = "  +4.27%  " }
    @{ Row = 35; D = $null;        E = "  +0.02%  " }
    @{ Row = 36; D = $null;        E = "  +1.19%  " }
    @{ Row = 37; D = "169.72";     E = "  -1.30%  " }
    @{ Row = 38; D = "6.94";       E = "  +1.22%  " }
    @{ Row = 39; D = "1.56";       E = "  +6.50%  " }
    @{ Row = 40; D = "5.02";       E = "  +5.79%  " }
    @{ Row = 41; D = "0.0808";     E = "  +4.27%  " }
    @{ Row = 42; D = "27.36";      E = "  +11.40%  " }
    @{ Row = 43; D = $null;        E = "  +1.54%  " }
    @{ Row = 44; D = $null;        E = "  +0.66%  " }
    @{ Row = 45; D = "0.999";      E = "  -0.14%  " }
    @{ Row = 46; D = $null;        E = "  +2.79%  " }
    @{ Row = 47; D = "1.21";       E = "  +5.42%  " }
    @{ Row = 48; D = "1.66";       E = "  +2.65%  " }
    @{ Row = 49; D = "2.505.23";   E = "  +12.11%  " }
    @{ Row = 50; D = "6.94";       E = "  +4.14%  " }
    @{ Row = 51; D = "2.40";       E = "  +11.18%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $text = $u.D
        # Guard against Excel's automatic "looks like a number" coercion
        # (e.g. "600.98" -> 600.98) so the cell keeps its original text
        # semantics, same as the source inline-string value. Values that
        # contain more than one '.' (e.g. "65.061.75") are already safe,
        # since Excel cannot parse them as a single numeric literal.
        $isNumeric = $text -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$'

        if ($isNumeric) {
            $ws.Range("D$row").Value = "'" + $text
        } else {
            $ws.Range("D$row").Value = $text
        }
    }

    $ws.Range("E$row").Value = $u.E
}
